$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.958.37"
$ws.Range("E2").Value = "  -1.74%  "

$ws.Range("D3").Value = "2.933.96"
$ws.Range("E3").Value = "  -2.28%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.36"
$ws.Range("E5").Value = "  -1.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.21"
$ws.Range("E6").Value = "  +0.81%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "2.914.52"
$ws.Range("E8").Value = "  -2.82%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.502"
$ws.Range("E9").Value = "  -2.92%  "

$ws.Range("E10").Value = "  +6.86%  "

$ws.Range("E11").Value = "  -2.32%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.450"
$ws.Range("E12").Value = "  -1.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000225"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.64"
$ws.Range("E14").Value = "  +1.02%  "

$ws.Range("E15").Value = "  -0.40%  "

$ws.Range("D16").Value = "3.417.84"
$ws.Range("E16").Value = "  -2.36%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.86"
$ws.Range("E17").Value = "  -1.19%  "

$ws.Range("D18").Value = "60.932.21"
$ws.Range("E18").Value = "  -1.62%  "

$ws.Range("D19").Value = "2.930.61"
$ws.Range("E19").Value = "  -2.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "429.25"
$ws.Range("E20").Value = "  -4.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.85"
$ws.Range("E21").Value = "  -1.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.674"
$ws.Range("E22").Value = "  -1.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.24"
$ws.Range("E23").Value = "  -1.69%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.74"
$ws.Range("E24").Value = "  -1.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.91"
$ws.Range("E25").Value = "  -1.40%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.18"
$ws.Range("E26").Value = "  -2.73%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.93"
$ws.Range("E27").Value = "  -1.92%  "

$ws.Range("E28").Value = "  +0.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.35"
$ws.Range("E29").Value = "  +1.75%  "

$ws.Range("E30").Value = "  +0.00%  "

$ws.Range("E31").Value = "  +4.91%  "

$ws.Range("E32").Value = "  -2.89%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.91"
$ws.Range("E33").Value = "  -1.32%  "

$ws.Range("E34").Value = "  -2.61%  "

$ws.Range("D35").Value = "0.0₃0832"
$ws.Range("E35").Value = "  -0.74%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.01"
$ws.Range("E36").Value = "  -1.42%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.71"
$ws.Range("E37").Value = "  -1.83%  "

$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.00"
$ws.Range("E38").Value = "  +3.13%  "

$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.90"
$ws.Range("E39").Value = "  -0.96%  "

$ws.Range("E40").Value = "  +1.81%  "

$ws.Range("E41").Value = "  -0.96%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.79"
$ws.Range("E42").Value = "  -2.30%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.290"
$ws.Range("E43").Value = "  +6.42%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.86"
$ws.Range("E44").Value = "  +3.09%  "

$ws.Range("E45").Value = "  -0.95%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "372.44"
$ws.Range("E46").Value = "  -6.73%  "

$ws.Range("D47").Value = "2.658.83"
$ws.Range("E47").Value = "  -2.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.32"
$ws.Range("E48").Value = "  +0.96%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.41"
$ws.Range("E49").Value = "  +7.51%  "

$ws.Range("B50").Value = "USDe"
$ws.Range("C50").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -0.01%  "

$ws.Range("E51").Value = "  -0.80%  "
